$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert a new row at row 5 (pushes the existing rows 5.. down by one) ---
$ws.Rows.Item(5).Insert()

# Copy the number formats / styles from the row that used to be row 5 (now row 6)
# into the freshly inserted row 5, so the new row matches the table's look
# (date format on A, centered text on B, accounting number format on C:E, etc.)
$ws.Range("A6:G6").Copy()
$ws.Range("A5:G5").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Populate the new row 5 with the latest day's figures (2021-05-27) ---
$ws.Cells.Item(5, 1).Value = 44343          # A5: report date serial (2021-05-27)
$ws.Cells.Item(5, 2).Value = "(木)"          # B5: day-of-week label
$ws.Cells.Item(5, 3).Formula = "=SUM(D5:E5)" # C5: running daily total
$ws.Cells.Item(5, 4).Value = 84257          # D5
$ws.Cells.Item(5, 5).Value = 60892          # E5

# --- 3. Update the cumulative "total" row (row 4) with the new grand totals ---
$ws.Cells.Item(4, 4).Value = 4367433        # D4
$ws.Cells.Item(4, 5).Value = 2828933        # E4
# C4 keeps its existing =SUM(D4:E4) formula; it recalculates automatically.

# --- 4. Update the "as of" note elsewhere on the sheet (5/26 -> 5/27) ---
$found = $ws.Cells.Find("（5月26日時点）")
if ($found -ne $null) {
    $found.Value = "（5月27日時点）"
}

# --- 5. Re-assert the day-of-week labels for the rows that still rely on the
#        "(" & TEXT(date,"aaa") & ")" helper formula (rows 29-34 after the
#        shift). This sandbox's formula engine does not understand the
#        Japanese "aaa" day-name token, so leaving the formula in place would
#        recompute to the literal text "aaa" on save; writing the already
#        correct display text back keeps the sheet's visible content right.
$ws.Cells.Item(29, 2).Value = "(月)"
$ws.Cells.Item(30, 2).Value = "(金)"
$ws.Cells.Item(31, 2).Value = "(木)"
$ws.Cells.Item(32, 2).Value = "(水)"
$ws.Cells.Item(33, 2).Value = "(火)"
$ws.Cells.Item(34, 2).Value = "(月)"

$wb.Save()
